$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 6 (pushes existing rows 6-11 down to 7-12),
# duplicating the "C3" column row (row 5) with an updated D/M value.
$ws.Rows.Item(6).Insert()

# Copy row 5 (C3) into the newly-inserted row 6
$ws.Range("A5:M5").Copy()
$ws.Range("A6:M6").PasteSpecial()

# Apply the data differences for the new row (D=5, M=10)
$ws.Range("D6").Value = 5
$ws.Range("M6").Value = 10

# Update the selection to match the saved state
$ws.Range("M7").Select() | Out-Null
